# Generate Report for Handback
# Updates the localization-status report: marks items as handed back (in
# sync with en-US), refreshes the "Latest Handback DateTime" timestamps,
# clears the stale "handback not latest" error message, and widens the
# Status / Latest Handback DateTime / Error Detail columns so the new,
# longer text fits.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-29 18:53:54"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-29 18:54:04"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
